$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 13 (log entry #6): fill in date, start/stop times, delta minutes,
# activity ("Summarising") and comment ("Writing  JavaScript  Notes")
$ws.Range("C13").Value = 43514
$ws.Range("D13").Value = 0.75347222222222221
$ws.Range("E13").Value = 0.98958333333333337
$ws.Range("G13").Value = 280
$ws.Range("H13").Value = "Summarising"
$ws.Range("I13").Value = "Writing  JavaScript  Notes"

# Update the active selection to match the saved view state
$ws.Range("H15").Select()
